# Updated: ut 13. 04. 2021
# Applies corrected AgTests (F) / AgPosit (G) figures for rows 362-403.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F362").Value = 228948
$ws.Range("F363").Value = 188364
$ws.Range("G363").Value = 2762
$ws.Range("F364").Value = 168082
$ws.Range("F365").Value = 184273
$ws.Range("F366").Value = 340030
$ws.Range("F367").Value = 766147
$ws.Range("G367").Value = 3918
$ws.Range("F368").Value = 343367
$ws.Range("G368").Value = 2273
$ws.Range("F369").Value = 233816
$ws.Range("F370").Value = 182095
$ws.Range("F371").Value = 159666
$ws.Range("G371").Value = 1954
$ws.Range("F372").Value = 179323
$ws.Range("G372").Value = 1861
$ws.Range("F373").Value = 349096
$ws.Range("F374").Value = 772546
$ws.Range("F375").Value = 348702
$ws.Range("G375").Value = 1846
$ws.Range("F376").Value = 220221
$ws.Range("F377").Value = 176836
$ws.Range("G377").Value = 1826
$ws.Range("F378").Value = 157180
$ws.Range("F379").Value = 179083
$ws.Range("F382").Value = 357155
$ws.Range("G382").Value = 1566
$ws.Range("F383").Value = 220772
$ws.Range("F385").Value = 150740
$ws.Range("F386").Value = 182225
$ws.Range("G386").Value = 1359
$ws.Range("F387").Value = 351532
$ws.Range("F388").Value = 728214
$ws.Range("F389").Value = 353389
$ws.Range("F390").Value = 220304
$ws.Range("F391").Value = 176481
$ws.Range("F392").Value = 218539
$ws.Range("G392").Value = 1206
$ws.Range("F393").Value = 298736
$ws.Range("G393").Value = 1189
$ws.Range("F394").Value = 162221
$ws.Range("G394").Value = 616
$ws.Range("F395").Value = 738233
$ws.Range("G395").Value = 1919
$ws.Range("F396").Value = 163677
$ws.Range("G396").Value = 547
$ws.Range("F397").Value = 106174
$ws.Range("G397").Value = 630
$ws.Range("F398").Value = 290950
$ws.Range("G398").Value = 1442
$ws.Range("F399").Value = 194730
$ws.Range("G399").Value = 956
$ws.Range("F400").Value = 144906
$ws.Range("F401").Value = 261544
$ws.Range("G401").Value = 909
$ws.Range("F402").Value = 688855
$ws.Range("G402").Value = 1334
$ws.Range("F403").Value = 330390
$ws.Range("G403").Value = 695
